$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are written as exact text (preserve trailing zeros / formatting)
$textCells = @("D5", "D6", "D8", "D10", "D11", "D12", "D13", "D15", "D19", "D20", "D21", "D22", "D23", "D25", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D36", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated Price (D) and Volume(1h) (E) values row by row
$ws.Range("D2").Value = "63.629.54"
$ws.Range("E2").Value = "  +4.76%  "
$ws.Range("D3").Value = "2.729.84"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "577.30"
$ws.Range("E5").Value = "  -2.21%  "
$ws.Range("D6").Value = "154.39"
$ws.Range("E6").Value = "  +6.55%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "0.610"
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("D9").Value = "2.755.78"
$ws.Range("E9").Value = "  +3.08%  "
$ws.Range("D10").Value = "6.71"
$ws.Range("E10").Value = "  +2.24%  "
$ws.Range("D11").Value = "0.113"
$ws.Range("E11").Value = "  +5.55%  "
$ws.Range("D12").Value = "0.162"
$ws.Range("E12").Value = "  +4.96%  "
$ws.Range("D13").Value = "0.389"
$ws.Range("E13").Value = "  +4.09%  "
$ws.Range("D14").Value = "3.221.44"
$ws.Range("E14").Value = "  +3.04%  "
$ws.Range("D15").Value = "26.39"
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("D16").Value = "63.567.85"
$ws.Range("E16").Value = "  +4.71%  "
$ws.Range("E17").Value = "  +6.27%  "
$ws.Range("D18").Value = "2.747.82"
$ws.Range("E18").Value = "  +3.46%  "
$ws.Range("D19").Value = "11.96"
$ws.Range("E19").Value = "  +3.71%  "
$ws.Range("D20").Value = "4.87"
$ws.Range("E20").Value = "  +2.87%  "
$ws.Range("D21").Value = "360.07"
$ws.Range("E21").Value = "  +2.84%  "
$ws.Range("D22").Value = "6.98"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").Value = "65.99"
$ws.Range("E25").Value = "  +3.16%  "
$ws.Range("E26").Value = "  +4.86%  "
$ws.Range("D27").Value = "8.54"
$ws.Range("E27").Value = "  +5.11%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("D29").Value = "0.0₃0909"
$ws.Range("E29").Value = "  +12.35%  "
$ws.Range("D30").Value = "1.99"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").Value = "7.13"
$ws.Range("E31").Value = "  +5.76%  "
$ws.Range("D32").Value = "172.57"
$ws.Range("E32").Value = "  +3.14%  "
$ws.Range("D33").Value = "1.23"
$ws.Range("E33").Value = "  +15.58%  "
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("E35").Value = "  +3.24%  "
$ws.Range("D36").Value = "4.77"
$ws.Range("E36").Value = "  +7.22%  "
$ws.Range("E37").Value = "  +9.05%  "
$ws.Range("D38").Value = "1.81"
$ws.Range("E38").Value = "  +9.89%  "
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +14.27%  "
$ws.Range("D40").Value = "343.31"
$ws.Range("E40").Value = "  +4.85%  "
$ws.Range("D41").Value = "4.22"
$ws.Range("E41").Value = "  +5.58%  "
$ws.Range("D42").Value = "39.38"
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("E43").Value = "  +6.45%  "
$ws.Range("D44").Value = "21.70"
$ws.Range("E44").Value = "  +5.52%  "
$ws.Range("D45").Value = "21.83"
$ws.Range("E45").Value = "  +6.50%  "
$ws.Range("D46").Value = "140.10"
$ws.Range("E46").Value = "  +4.09%  "
$ws.Range("D47").Value = "0.0591"
$ws.Range("E47").Value = "  +6.07%  "
$ws.Range("D48").Value = "0.645"
$ws.Range("E48").Value = "  +5.08%  "
$ws.Range("D49").Value = "0.0254"
$ws.Range("E49").Value = "  +3.32%  "
$ws.Range("E50").Value = "  +1.48%  "
$ws.Range("D51").Value = "0.996"
$ws.Range("E51").Value = "  -0.37%  "
